$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to Text format so numeric-looking strings
# (e.g. "238.42", "5.220") are preserved exactly as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.152.57'
$ws.Range("E2").Value = '  -3.21%  '
$ws.Range("D3").Value = '1.848.57'
$ws.Range("E3").Value = '  -2.24%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '0.7022'
$ws.Range("E5").Value = '  -5.12%  '
$ws.Range("D6").Value = '238.42'
$ws.Range("E6").Value = '  -1.76%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '0.3052'
$ws.Range("E8").Value = '  -3.93%  '
$ws.Range("D9").Value = '0.07467'
$ws.Range("E9").Value = '  +4.12%  '
$ws.Range("D10").Value = '23.38'
$ws.Range("E10").Value = '  -6.06%  '
$ws.Range("D11").Value = '0.08127'
$ws.Range("E11").Value = '  -2.54%  '
$ws.Range("D12").Value = '1.878.27'
$ws.Range("E12").Value = '  -0.94%  '
$ws.Range("D13").Value = '0.7252'
$ws.Range("E13").Value = '  -4.47%  '
$ws.Range("D14").Value = '5.220'
$ws.Range("E14").Value = '  -3.75%  '
$ws.Range("D15").Value = '88.68'
$ws.Range("E15").Value = '  -4.80%  '
$ws.Range("D16").Value = '29.266.89'
$ws.Range("E16").Value = '  -2.93%  '
$ws.Range("D17").Value = '5.758'
$ws.Range("E17").Value = '  -6.56%  '
$ws.Range("D18").Value = '238.16'
$ws.Range("E18").Value = '  -5.31%  '
$ws.Range("D19").Value = '13.06'
$ws.Range("E19").Value = '  -4.15%  '
$ws.Range("D20").Value = '0.000007624'
$ws.Range("E20").Value = '  -3.14%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").Value = '2.118.60'
$ws.Range("E22").Value = '  -3.79%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '7.578'
$ws.Range("E24").Value = '  -4.81%  '
$ws.Range("D25").Value = '8.989'
$ws.Range("E25").Value = '  -3.52%  '
$ws.Range("D26").Value = '161.27'
$ws.Range("E26").Value = '  -2.09%  '
$ws.Range("D27").Value = '0.1451'
$ws.Range("E27").Value = '  -7.68%  '
$ws.Range("D28").Value = '18.06'
$ws.Range("E28").Value = '  -3.65%  '
$ws.Range("D29").Value = '1.983'
$ws.Range("E29").Value = '  -3.81%  '
$ws.Range("E30").Value = '  -5.54%  '
$ws.Range("D31").Value = '4.547'
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("E32").Value = '  -2.88%  '
$ws.Range("D33").Value = '3.971'
$ws.Range("E33").Value = '  -5.55%  '
$ws.Range("D34").Value = '0.05174'
$ws.Range("E34").Value = '  -3.44%  '
$ws.Range("D35").Value = '1.186'
$ws.Range("E35").Value = '  -5.64%  '
$ws.Range("D36").Value = '1.038'
$ws.Range("E36").Value = '  +3.44%  '
$ws.Range("D37").Value = '0.6998'
$ws.Range("E37").Value = '  -9.59%  '
$ws.Range("D38").Value = '2.656'
$ws.Range("E38").Value = '  -2.54%  '
$ws.Range("D39").Value = '0.01862'
$ws.Range("E39").Value = '  -5.00%  '
$ws.Range("D40").Value = '2.677'
$ws.Range("E40").Value = '  -3.09%  '
$ws.Range("D41").Value = '0.9339'
$ws.Range("E41").Value = '  +6.65%  '
$ws.Range("D42").Value = '6.018'
$ws.Range("E42").Value = '  -1.07%  '
$ws.Range("D43").Value = '1.077.53'
$ws.Range("E43").Value = '  -2.36%  '
$ws.Range("D44").Value = '0.4283'
$ws.Range("E44").Value = '  -6.27%  '
$ws.Range("D45").Value = '70.29'
$ws.Range("E45").Value = '  -3.17%  '
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").Value = '102.53'
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("D48").Value = '2.008.08'
$ws.Range("E48").Value = '  -2.82%  '
$ws.Range("D49").Value = '1.743'
$ws.Range("E49").Value = '  -6.59%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.165'
$ws.Range("E50").Value = '  -4.66%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '7.031'
$ws.Range("E51").Value = '  -7.42%  '

# Restore default formatting/style on column D so style indices match the source
# (avoids leaving a stray custom number format behind).
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
